# "save data done + era data updated"
# Add a new "Save" column (H) to the stats sheet: a header cell matching the
# existing header formatting, and zero-valued data cells for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", formatted like the other header cells (B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New data column H2:H8, all zeros for now.
$ws.Range("H2:H8").Value = 0
